$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 1315.75
$ws.Cells.Item(135, 9).Value = 1046.5952
$ws.Cells.Item(135, 10).Value = 6968
$ws.Cells.Item(135, 11).Value = 9419.3568
$ws.Cells.Item(135, 12).Value = 62712
$ws.Cells.Item(135, 13).Value = -6884.3568
$ws.Cells.Item(135, 14).Value = -67782
$ws.Cells.Item(137, 8).Value = 17858088
$ws.Cells.Item(137, 9).Value = 21739748
$ws.Cells.Item(137, 10).Value = 2451.8
$ws.Cells.Item(137, 11).Value = 65219244
$ws.Cells.Item(137, 12).Value = 7355.400000000001
$ws.Cells.Item(137, 13).Value = -65216694
$ws.Cells.Item(137, 14).Value = -12455.4
$ws.Cells.Item(138, 8).Value = 1079.4343
$ws.Cells.Item(138, 9).Value = 479.125
$ws.Cells.Item(138, 10).Value = 1644.4314
$ws.Cells.Item(138, 11).Value = 1437.375
$ws.Cells.Item(138, 12).Value = 4933.2942
$ws.Cells.Item(138, 13).Value = 3702.625
$ws.Cells.Item(138, 14).Value = -15213.2942
$ws.Cells.Item(141, 8).Value = 2007.7354
$ws.Cells.Item(141, 9).Value = 1124.614
$ws.Cells.Item(141, 10).Value = 6583.909
$ws.Cells.Item(141, 11).Value = 3373.842
$ws.Cells.Item(141, 12).Value = 19751.727
$ws.Cells.Item(141, 13).Value = 1806.158
$ws.Cells.Item(141, 14).Value = -30111.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 23473.207
$ws.Cells.Item(32, 9).Value = 4221.32
$ws.Cells.Item(32, 10).Value = 344338
$ws.Cells.Item(32, 11).Value = 4221.32
$ws.Cells.Item(32, 12).Value = 344338
$ws.Cells.Item(32, 13).Value = -3934.32
$ws.Cells.Item(32, 14).Value = -344912
$ws.Cells.Item(63, 8).Value = 32897.25
$ws.Cells.Item(63, 9).Value = 32897.25
$ws.Cells.Item(63, 11).Value = 32897.25
$ws.Cells.Item(63, 13).Value = -32211.25
$ws.Cells.Item(66, 8).Value = 32897.25
$ws.Cells.Item(66, 9).Value = 32897.25
$ws.Cells.Item(66, 11).Value = 164486.25
$ws.Cells.Item(66, 13).Value = -161054.25
$ws.Cells.Item(74, 8).Value = 4213.8096
$ws.Cells.Item(74, 9).Value = 1089.4546
$ws.Cells.Item(74, 10).Value = 15669.777
$ws.Cells.Item(74, 11).Value = 1089.4546
$ws.Cells.Item(74, 12).Value = 15669.777
$ws.Cells.Item(74, 13).Value = -215.4546
$ws.Cells.Item(74, 14).Value = -17417.777
$ws.Cells.Item(77, 8).Value = 4213.8096
$ws.Cells.Item(77, 9).Value = 1089.4546
$ws.Cells.Item(77, 10).Value = 15669.777
$ws.Cells.Item(77, 11).Value = 5447.273
$ws.Cells.Item(77, 12).Value = 78348.88499999999
$ws.Cells.Item(77, 13).Value = -1079.273
$ws.Cells.Item(77, 14).Value = -87084.88499999999
$ws.Cells.Item(97, 8).Value = 5054.1816
$ws.Cells.Item(97, 9).Value = 6060.0557
$ws.Cells.Item(97, 11).Value = 6060.0557
$ws.Cells.Item(97, 13).Value = -5564.0557
$ws.Cells.Item(102, 8).Value = 1474.75
$ws.Cells.Item(102, 9).Value = 950
$ws.Cells.Item(102, 10).Value = 1999.5
$ws.Cells.Item(102, 11).Value = 950
$ws.Cells.Item(102, 12).Value = 1999.5
$ws.Cells.Item(102, 13).Value = 672
$ws.Cells.Item(102, 14).Value = -5243.5
$ws.Cells.Item(122, 8).Value = 1388.8788
$ws.Cells.Item(122, 9).Value = 1258.4231
$ws.Cells.Item(122, 11).Value = 3775.2693
$ws.Cells.Item(122, 13).Value = -1325.2693

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1209.9445
$ws.Cells.Item(20, 9).Value = 1071.2
$ws.Cells.Item(20, 11).Value = 1071.2
$ws.Cells.Item(20, 13).Value = -824.2
$ws.Cells.Item(134, 8).Value = 15153615
$ws.Cells.Item(134, 9).Value = 19232550
$ws.Cells.Item(134, 10).Value = 3290.5715
$ws.Cells.Item(134, 11).Value = 57697650
$ws.Cells.Item(134, 12).Value = 9871.7145
$ws.Cells.Item(134, 13).Value = -57695115
$ws.Cells.Item(134, 14).Value = -14941.7145
$ws.Cells.Item(135, 8).Value = 41039.8
$ws.Cells.Item(135, 10).Value = 41039.8
$ws.Cells.Item(135, 12).Value = 41039.8
$ws.Cells.Item(135, 14).Value = -51179.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1508.4062
$ws.Cells.Item(31, 9).Value = 935.5848999999999
$ws.Cells.Item(31, 10).Value = 4268.364
$ws.Cells.Item(31, 11).Value = 935.5848999999999
$ws.Cells.Item(31, 12).Value = 4268.364
$ws.Cells.Item(31, 13).Value = -640.5848999999999
$ws.Cells.Item(31, 14).Value = -4858.364
$ws.Cells.Item(34, 8).Value = 1508.4062
$ws.Cells.Item(34, 9).Value = 935.5848999999999
$ws.Cells.Item(34, 10).Value = 4268.364
$ws.Cells.Item(34, 11).Value = 935.5848999999999
$ws.Cells.Item(34, 12).Value = 4268.364
$ws.Cells.Item(34, 13).Value = -733.5848999999999
$ws.Cells.Item(34, 14).Value = -4672.364
$ws.Cells.Item(134, 8).Value = 2192.2246
$ws.Cells.Item(134, 9).Value = 1471.8889
$ws.Cells.Item(134, 11).Value = 4415.6667
$ws.Cells.Item(134, 13).Value = -1880.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(118, 8).Value = 3026.1428
$ws.Cells.Item(118, 9).Value = 964.5
$ws.Cells.Item(118, 10).Value = 3850.8
$ws.Cells.Item(118, 11).Value = 2893.5
$ws.Cells.Item(118, 12).Value = 11552.4
$ws.Cells.Item(118, 13).Value = -1650.5
$ws.Cells.Item(118, 14).Value = -14038.4
$ws.Cells.Item(131, 8).Value = 7938397.5
$ws.Cells.Item(131, 9).Value = 591.125
$ws.Cells.Item(131, 10).Value = 9806117
$ws.Cells.Item(131, 11).Value = 1773.375
$ws.Cells.Item(131, 12).Value = 29418351
$ws.Cells.Item(131, 13).Value = 3266.625
$ws.Cells.Item(131, 14).Value = -29428431
$ws.Cells.Item(140, 8).Value = 7057.7837
$ws.Cells.Item(140, 9).Value = 9109.52
$ws.Cells.Item(140, 10).Value = 2783.3333
$ws.Cells.Item(140, 11).Value = 27328.56
$ws.Cells.Item(140, 12).Value = 8349.999899999999
$ws.Cells.Item(140, 13).Value = -22148.56
$ws.Cells.Item(140, 14).Value = -18709.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 8000
$ws.Cells.Item(52, 10).Value = 8000
$ws.Cells.Item(52, 12).Value = 8000
$ws.Cells.Item(52, 14).Value = -8518
$ws.Cells.Item(102, 8).Value = 2120.4827
$ws.Cells.Item(102, 9).Value = 2245.7896
$ws.Cells.Item(102, 11).Value = 2245.7896
$ws.Cells.Item(102, 13).Value = -623.7896000000001
$ws.Cells.Item(108, 8).Value = 31900
$ws.Cells.Item(108, 10).Value = 31900
$ws.Cells.Item(108, 12).Value = 31900
$ws.Cells.Item(108, 14).Value = -39580
$ws.Cells.Item(132, 8).Value = 2548.4324
$ws.Cells.Item(132, 9).Value = 2231.2068
$ws.Cells.Item(132, 11).Value = 6693.6204
$ws.Cells.Item(132, 13).Value = -4163.6204

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 1984.6154
$ws.Cells.Item(68, 9).Value = 1537.5
$ws.Cells.Item(68, 10).Value = 2700
$ws.Cells.Item(68, 11).Value = 1537.5
$ws.Cells.Item(68, 12).Value = 2700
$ws.Cells.Item(68, 13).Value = -788.5
$ws.Cells.Item(68, 14).Value = -4198
$ws.Cells.Item(71, 8).Value = 1984.6154
$ws.Cells.Item(71, 9).Value = 1537.5
$ws.Cells.Item(71, 10).Value = 2700
$ws.Cells.Item(71, 11).Value = 7687.5
$ws.Cells.Item(71, 12).Value = 13500
$ws.Cells.Item(71, 13).Value = -3943.5
$ws.Cells.Item(71, 14).Value = -20988
$ws.Cells.Item(137, 8).Value = 35000
$ws.Cells.Item(137, 10).Value = 35000
$ws.Cells.Item(137, 12).Value = 35000
$ws.Cells.Item(137, 14).Value = -45200

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(6, 8).Value = 278025.25
$ws.Cells.Item(6, 10).Value = 506050.5
$ws.Cells.Item(6, 12).Value = 506050.5
$ws.Cells.Item(6, 14).Value = -506280.5
$ws.Cells.Item(45, 8).Value = 12385
$ws.Cells.Item(45, 10).Value = 12385
$ws.Cells.Item(45, 12).Value = 12385
$ws.Cells.Item(45, 14).Value = -13367
$ws.Cells.Item(136, 8).Value = 36075.586
$ws.Cells.Item(136, 9).Value = 56562.055
$ws.Cells.Item(136, 10).Value = 2552.2727
$ws.Cells.Item(136, 11).Value = 169686.165
$ws.Cells.Item(136, 12).Value = 7656.8181
$ws.Cells.Item(136, 13).Value = -167136.165
$ws.Cells.Item(136, 14).Value = -12756.8181
